$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("concept:name", "str"),
    @("case:concept:name", "str"),
    @("lifecycle:state", "str"),
    @("unsatisfied_condition_description", "str"),
    @("response_status_code", "float"),
    @("time:timestamp", "datetime"),
    @("identifier:id", "str"),
    @("planned_operation_time", "str"),
    @("complete_service_time", "str"),
    @("event_id", "str"),
    @("SubProcessID", "str"),
    @("process_model_id", "str"),
    @("org:resource", "str"),
    @("operation_end_time", "datetime"),
    @("lifecycle:transition", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("requested_service_url", "str"),
    @("case", "str"),
    @("current_task", "str"),
    @("parameters", "dict")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
